# household_new.xlsx edit:
# Insert two new columns (auxillaryHash / auxillaryHash.cell_type) into the
# "survey" sheet between "selectionArgs.cell_type" and "comments", and
# rewrite the "comments" example text in row 7 to document the new
# auxillaryHash setting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# 1. Insert two blank columns at J:K - this pushes the old J:N (comments,
#    line_text.elementName, line_text.isInstanceMetadata,
#    line_subtext.elementName, line_subtext.isInstanceMetadata) to L:P.
$ws.Range("J1:K1").EntireColumn.Insert()

# 2. Header row (row 1) for the two new columns.
$ws.Range("J1").Value = "auxillaryHash"
$ws.Range("K1").Value = "auxillaryHash.cell_type"

# 3. Data row (row 7) values for the two new columns.
#    J7 looks like a formula (leading apostrophe) so it needs to be entered
#    with a doubled leading apostrophe to keep one literal quote character
#    in the stored string while still tagging the cell with quote-prefix.
$ws.Range("J7").Value = "''household_id='+escape(data('household_id'))"
$ws.Range("K7").Value = "formula"

# 4. Column widths for the two new columns (existing J:N widths already
#    shifted automatically with the insert above).
$ws.Range("J1").ColumnWidth = 27.072916666666668
$ws.Range("K1").ColumnWidth = 19.346354166666668

# 5. Rewrite the comments cell (old J7, now L7) with the new description
#    that documents auxillaryHash ahead of the pre-existing
#    joined_through_name paragraph.
$commentsCell = $ws.Range("L7")
$commentsCell.Value = "auxillaryHash defines the auxillary hash to supply when creating a new sub-form. This is an ampersand-separated list of elementName=value pairs that will be used to initialize the subform. The joined_through_name value identifies the name (elementName) in the model that should be used when scanning in the joins lists for the table_id to discover the foreign key column to filter on in the subform.`nIf this is omitted, we would probably just scan the entire model to see if table_id appears anywhere and use the first match we find. "

$commentsCell.Characters(1, 13).Font.Bold = $true
$commentsCell.Characters(14, 180).Font.Bold = $false
$commentsCell.Characters(194, 19).Font.Bold = $true

# 6. Row 7 is taller to fit the longer comments text.
$ws.Range("A7").RowHeight = 220.5

# 7. Update the view so the newly inserted columns are visible / selected,
#    matching the author's saved cursor position.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("J8").Select()
